$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - South Korea
$ws.Range("A2").Value = "South Korea"
$ws.Range("C2").Value = 44012
$ws.Range("D2").Value = 0.02203125
$ws.Range("E2").Value = 0.02438820524453683
$ws.Range("F2").Value = 0.01820958489448707
$ws.Range("G2").Value = 0.00617862035004976
$ws.Range("H2").Value = 0.7466553898453096
$ws.Range("I2").Value = 0.2533446101546905

# Row 3 - China
$ws.Range("D3").Value = 0.02290248925501433
$ws.Range("E3").Value = 0.0235169659895225
$ws.Range("F3").Value = 0.01109164660972722
$ws.Range("G3").Value = 0.01242531937979528
$ws.Range("H3").Value = 0.4716444550997299
$ws.Range("I3").Value = 0.52835554490027

# Row 4 - Germany
$ws.Range("C4").Value = 44012
$ws.Range("D4").Value = 0.04641945524453683

# Row 5 - USA / All
$ws.Range("C5").Value = 44009
$ws.Range("D5").Value = 0.04752702796222197
$ws.Range("E5").Value = -0.001107572717685142
$ws.Range("F5").Value = 0.007542255352947743
$ws.Range("G5").Value = -0.008649828070632882
$ws.Range("H5").Value = 0.4657989435728767
$ws.Range("I5").Value = 0.5342010564271233

# Row 6 - USA / NYC
$ws.Range("C6").Value = 44012
$ws.Range("D6").Value = 0.08719684220304529
$ws.Range("E6").Value = -0.04077738695850846
$ws.Range("F6").Value = 0.007840018472217361
$ws.Range("G6").Value = -0.04861740543072583
$ws.Range("H6").Value = 0.1388660326708363
$ws.Range("I6").Value = 0.8611339673291638

# Row 7 - Spain
$ws.Range("C7").Value = 43972
$ws.Range("D7").Value = 0.121913536873179
$ws.Range("E7").Value = -0.07549408162864213
$ws.Range("F7").Value = -0.05133669685010826
$ws.Range("G7").Value = -0.02415738477853388
$ws.Range("H7").Value = 0.6800095549560449
$ws.Range("I7").Value = 0.3199904450439552

# Row 8 - Italy
$ws.Range("C8").Value = 44012
$ws.Range("D8").Value = 0.1403006799609075
$ws.Range("E8").Value = -0.09388122471637062
$ws.Range("F8").Value = -0.05876468611030111
$ws.Range("G8").Value = -0.0351165386060695
$ws.Range("H8").Value = 0.6259471612970338
$ws.Range("I8").Value = 0.3740528387029662
